$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 219 (pushes the existing rows 219..256 down to 220..257)
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Cells.Item(219, 1).Value = 3
$ws.Cells.Item(219, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(219, 3).Value = "Coquimbo"
$ws.Cells.Item(219, 4).Value = 44522
$ws.Cells.Item(219, 5).Value = 5
$ws.Cells.Item(219, 6).Value = 100112003
$ws.Cells.Item(219, 7).Value = "Ajo"
$ws.Cells.Item(219, 8).Value = "Chino"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 73
$ws.Cells.Item(219, 11).Value = 16000
$ws.Cells.Item(219, 12).Value = 16500
$ws.Cells.Item(219, 13).Value = 16260
$ws.Cells.Item(219, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(219, 15).Value = "China"
$ws.Cells.Item(219, 16).Value = 1626
$ws.Cells.Item(219, 17).Value = 10
$ws.Cells.Item(219, 18).Value = "Hortaliza"
